$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Real" week-of-completion values for the rows that were left blank.
$ws.Range("G27").Value = 3
$ws.Range("G28").Value = 3
$ws.Range("G30").Value = 3

# Move the current selection to match where the user ended up (G31).
$ws.Range("G31").Select()
